# Apply the edit described by the diff to LOM3049.xlsx
# Summary of change:
#  - The standalone "Docentes responsaveis" value row (old row 13, containing only
#    "5840521 - Rosa Ana Conte" in B/C with no label in A) is removed, shifting all
#    subsequent rows up by one.
#  - Several label/value cells are then updated to their new (shifted) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old standalone row 13 (B13/C13 = "5840521 - Rosa Ana Conte", no A13 label).
# This shifts every row below it up by one, matching the new dimension A1:C24.
$ws.Rows.Item(13).Delete()

# --- Update cell values to match the post-edit content ---

# Row 10 (Objetivos:) now holds the docente value instead of the long objectives text.
$ws.Range("B10").Value = '5840521 - Rosa Ana Conte'
$ws.Range("C10").Value = '5840521 - Rosa Ana Conte'

# Row 13 (Programa resumido:) now holds "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now holds the activation date.
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

# Row 18 (Metodo:) now holds the docente value.
$ws.Range("B18").Value = '5840521 - Rosa Ana Conte'
$ws.Range("C18").Value = '5840521 - Rosa Ana Conte'

# Row 19 (Criterio:) now holds the evaluation-method text.
$ws.Range("B19").Value = 'Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações.'
$ws.Range("C19").Value = 'Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações.'

# Row 20 (Norma de recuperacao:) now holds the criteria text.
$ws.Range("B20").Value = 'Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3.'
$ws.Range("C20").Value = 'Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3.'

# Row 21 (Bibliografia:) now holds the recovery-norm text.
$ws.Range("B21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
